# Update the "dSF" (column F) values for a handful of rows in Sheet1.
# These rows represent a repull/recalculation of the dSF data (per the
# commit message: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    5  = -1
    9  = 4
    10 = 2
    11 = -2
    13 = 3
    17 = 1
    21 = -3
    23 = 5
    28 = 1
    33 = 3
    38 = -1
    42 = 7
    46 = -6
    53 = 1
    55 = 3
    57 = -5
    58 = 0
    61 = 1
    64 = -2
    69 = -2
    70 = -5
    72 = 1
    75 = 1
    77 = 4
    80 = 0
    81 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
